$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Prefix the "Bugfix" task rows (4-13 in column A) with "Fix: " (a couple of
# them already started with "Fix" in a slightly different form, so set the
# full target string explicitly to match exactly). Edited in the same
# order as the original authoring session so newly-minted shared strings
# land in the same slots.
$ws.Range("A7").Value = "Fix: func_train in subway fucked up"
$ws.Range("A8").Value = "Fix: Fucked up chair collision"
$ws.Range("A6").Value = "Fix: relative velocity being added multiple times"
$ws.Range("A5").Value = "Fix: game connecting to server during intro"
$ws.Range("A4").Value = 'Fix: "velocity increasing but not speed" (prediction miss bug when hitting wall at specific angle sometimes)'
$ws.Range("A9").Value = "Fix: Zombification: you can get stuck"
$ws.Range("A13").Value = "Fix: Incorrect update channel description (Release) being pulled on Debug"
$ws.Range("A12").Value = "Fix: Alt+tab crash in dedicated fullscreen"
$ws.Range("A11").Value = "Fix: Game does not clear entities on 2nd entry into a map"
$ws.Range("A10").Value = "Fix:  Green pipe, add poster there"

# Move the active selection like the author's session ended up (A23).
$ws.Range("A23").Select()
